$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Create the new hidden "Tabelle2" sheet (boolean list values) right after Tabelle1.
# Writing these values first makes them land at shared-string indices 8..10,
# matching the order produced by the original authoring session.
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "Tabelle2"
$ws2.PageSetup.TopMargin = 0.78740157499999996 * 72
$ws2.PageSetup.BottomMargin = 0.78740157499999996 * 72
$ws2.Range("A1").Value = "boolean"
$ws2.Range("A2").Value = "ja"
$ws2.Range("A3").Value = "nein"
$ws2.Range("A1:A3").Select() | Out-Null
$ws2.Visible = $false

# Go back to Tabelle1 and insert the new "Juenger als 23..." row at row 3,
# pushing the existing rows 3-6 down to 4-7.
$ws1.Select() | Out-Null
$ws1.Rows.Item(3).Insert()
$ws1.Range("A3").Value = "Juenger als 23 oder geboren vor 1940"
$ws1.Range("B3").Value = "nein"

# Restrict B3 to a dropdown list sourced from Tabelle2!$A$2:$A$3 ("ja"/"nein").
$validation = $ws1.Range("B3").Validation
$validation.Add(3, 1, 1, "=Tabelle2!`$A`$2:`$A`$3")
$validation.IgnoreBlank = $true
$validation.InCellDropdown = $true
$validation.ShowInput = $true
$validation.ShowError = $true

# Match the author's final selection (cell A3 highlighted on Tabelle1).
$ws1.Range("A3").Select() | Out-Null
